$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A ("ID") and shift the rest right.
$ws.Range("A1").EntireColumn.Insert()

# Header row
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Matcher"
$ws.Range("C1").Value = "Target URL"
$ws.Range("D1").Value = "Type"
$ws.Range("E1").Value = "Info"
$ws.Range("F1").Value = "Auto Redirect"

# Row 2
$ws.Range("A2").Formula = '=""'
$ws.Range("D2").Value = "partial"
$ws.Range("E2").Value = "Sample Redirect"
$ws.Range("F2").Value = $false

# Row 3
$ws.Range("A3").Formula = '=""'
$ws.Range("D3").Value = "wildcard"
$ws.Range("E3").Value = "Legacy Section Redirect"
$ws.Range("F3").Value = $true
